$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.941.53'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  +1.88%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.639.91'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  +3.92%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.89'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '199.94'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +2.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.221'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +9.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.648'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.18'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +1.48%  '
$ws.Range('E12').Value = '  +2.00%  '
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.215.34'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +3.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '667.89'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  +12.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.03'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  +1.99%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '71.073.65'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  +1.81%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.638.27'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +3.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.12'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('E21').Value = '  +1.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.51'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +1.17%  '
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '105.57'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +3.51%  '
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.01'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -4.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.50'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -3.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.82'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +2.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.19'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +2.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.65'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +8.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.25'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +3.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.24'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  -1.23%  '
$ws.Range('E33').Value = '  +0.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.50'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +0.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.989.64'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  +7.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0₃0877'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +6.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.04'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -1.61%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('B39').NumberFormat = 'General'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('C39').NumberFormat = 'General'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.02'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  +1.86%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('B40').NumberFormat = 'General'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('C40').NumberFormat = 'General'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '506.74'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  +5.94%  '
$ws.Range('E41').Value = '  -0.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.54'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -2.97%  '
$ws.Range('E43').Value = '  +2.00%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('B44').NumberFormat = 'General'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('C44').NumberFormat = 'General'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.09'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +9.55%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'VeChain'
$ws.Range('B45').NumberFormat = 'General'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C45').NumberFormat = 'General'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0461'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +1.86%  '
$ws.Range('E46').Value = '  +6.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.141'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.70'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +3.50%  '
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000249'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +1.66%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'LidoDAOToken'
$ws.Range('B51').NumberFormat = 'General'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C51').NumberFormat = 'General'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.96'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  +5.26%  '
